$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 182 (13-09-2021)
$ws.Range("B182").Value = 34869.6
$ws.Range("C182").Value = 15105.6
$ws.Range("D182").Value = 7068.4
$ws.Range("E182").Value = 30447.4
$ws.Range("F182").Value = 6676.9
$ws.Range("G182").Value = 15701.4
$ws.Range("H182").Value = 3127.9
$ws.Range("I182").Value = 4991.7
$ws.Range("J182").Value = 1570.1
$ws.Range("K182").Value = 71804
$ws.Range("L182").Value = 4041.2
$ws.Range("M182").Value = 1633.8
$ws.Range("N182").Value = 17446.3
$ws.Range("O182").Value = 1438.3
$ws.Range("P182").Value = 80065.9
$ws.Range("Q182").Value = 116403.7
$ws.Range("R182").Value = 1312.4
$ws.Range("S182").Value = 51496.7
$ws.Range("T182").Value = 18063.7

# Row 183 (14-09-2021)
$ws.Range("A183").Value = "14-09-2021"
$ws.Range("B183").Value = 34577.6
$ws.Range("C183").Value = 15037.8
$ws.Range("D183").Value = 7034.1
$ws.Range("E183").Value = 30670.1
$ws.Range("F183").Value = 6653
$ws.Range("G183").Value = 15723
$ws.Range("H183").Value = 3148.8
$ws.Range("I183").Value = 4917.2
$ws.Range("J183").Value = 1555.5
$ws.Range("K183").Value = 72040.7
$ws.Range("L183").Value = 4055.7
$ws.Range("M183").Value = 1623.8
$ws.Range("N183").Value = 17434.9
$ws.Range("O183").Value = 1434.4
$ws.Range("P183").Value = 78989.2
$ws.Range("Q183").Value = 116180.6
$ws.Range("R183").Value = 1310.2
$ws.Range("S183").Value = 51815.2
$ws.Range("T183").Value = 17895.9

# Row 184 (15-09-2021)
$ws.Range("A184").Value = "15-09-2021"
$ws.Range("B184").Value = 34814.4
$ws.Range("C184").Value = 15161.5
$ws.Range("D184").Value = 7016.5
$ws.Range("E184").Value = 30511.7
$ws.Range("F184").Value = 6583.6
$ws.Range("G184").Value = 15616
$ws.Range("H184").Value = 3153.4
$ws.Range("I184").Value = 4867.3
$ws.Range("J184").Value = 1555.3
$ws.Range("K184").Value = 71711.5
$ws.Range("L184").Value = 4065.5
$ws.Range("M184").Value = 1628
$ws.Range("N184").Value = 17354
$ws.Range("O184").Value = 1424
$ws.Range("P184").Value = 78488.5
$ws.Range("Q184").Value = 115062.5
$ws.Range("R184").Value = 1322.1
$ws.Range("S184").Value = 52192.3
$ws.Range("T184").Value = 18008.8

# Row 185 (16-09-2021)
$ws.Range("A185").Value = "16-09-2021"
$ws.Range("B185").Value = 34751.3
$ws.Range("C185").Value = 15181.9
$ws.Range("D185").Value = 7027.5
$ws.Range("E185").Value = 30323.3
$ws.Range("F185").Value = 6622.6
$ws.Range("G185").Value = 15651.8
$ws.Range("H185").Value = 3130.1
$ws.Range("I185").Value = 4807.7
$ws.Range("K185").Value = 71674.7
$ws.Range("L185").Value = 4045.1
$ws.Range("M185").Value = 1631.7
$ws.Range("N185").Value = 17278.7
$ws.Range("O185").Value = 1418
$ws.Range("P185").Value = 79738.1
$ws.Range("Q185").Value = 113794.3
$ws.Range("R185").Value = 1323.4
$ws.Range("T185").Value = 17885.1

# Row 186 (17-09-2021)
$ws.Range("A186").Value = "17-09-2021"
$ws.Range("B186").Value = 34584.9
$ws.Range("C186").Value = 15044
$ws.Range("D186").Value = 6963.6
$ws.Range("E186").Value = 30500.1
$ws.Range("F186").Value = 6570.2
$ws.Range("G186").Value = 15490.2
$ws.Range("H186").Value = 3140.5
$ws.Range("I186").Value = 4855.9
$ws.Range("J186").Value = 1548.5
$ws.Range("K186").Value = 71107.6
$ws.Range("L186").Value = 4035.2
$ws.Range("M186").Value = 1625.7
$ws.Range("N186").Value = 17276.8
$ws.Range("O186").Value = 1419.4
$ws.Range("P186").Value = 78512.1
$ws.Range("Q186").Value = 111439.4
$ws.Range("R186").Value = 1321.9
$ws.Range("S186").Value = 51307.7
$ws.Range("T186").Value = 17979.7

# Row 187 (20-09-2021)
$ws.Range("A187").Value = "20-09-2021"
$ws.Range("D187").Value = 6834.7
$ws.Range("F187").Value = 6394.2
$ws.Range("G187").Value = 15047.3
$ws.Range("J187").Value = 1527.9
$ws.Range("K187").Value = 69021.3
$ws.Range("L187").Value = 3958.9
$ws.Range("M187").Value = 1603.1
$ws.Range("O187").Value = 1387.1
